$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (row 1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update data row 2
$ws.Range("B2").Value = 137.51708750675024
$ws.Range("C2").Value = 178.301190856668
$ws.Range("D2").Value = 135.29728541516977
$ws.Range("E2").Value = 178.80066326742468

# Update data row 3
$ws.Range("B3").Value = 128.60867508149147
$ws.Range("C3").Value = 177.76980926461107
$ws.Range("D3").Value = 132.45973953741387
$ws.Range("E3").Value = 174.70118821497618

# Update the selection to match the new range
$ws.Range("B1:E3").Select()
